{"js": "// Replace each paragraph's text in document order with the target value.\n// The document has 1 date paragraph followed by 100 table-cell paragraphs\n// (20 rows x 5 cols); this edit only changes the text content of each run,\n// not the table/paragraph structure, so we can safely walk paragraphs in\n// order and do an in-place text replace on each one.\nconst newTexts = [\"2025-03-21 Friday\", \"78-78=\", \"45-1=\", \"4+38=\", \"3+83=\", \"83-42=\", \"72-36=\", \"54+30=\", \"37-13=\", \"60-29=\", \"46-13=\", \"82-1=\", \"91-65=\", \"57-22=\", \"75-7=\", \"69+10=\", \"42+0=\", \"16-12=\", \"39+2=\", \"78-49=\", \"77-52=\", \"72+10=\", \"47-3=\", \"9+6=\", \"31+22=\", \"63+6=\", \"60+16=\", \"77-25=\", \"78-47=\", \"30+28=\", \"24+61=\", \"32-20=\", \"28-19=\", \"89-64=\", \"57-34=\", \"14+72=\", \"23+53=\", \"17+16=\", \"50+45=\", \"30+46=\", \"17+56=\", \"42+15=\", \"13+32=\", \"74-37=\", \"20+57=\", \"55+39=\", \"45+17=\", \"55+19=\", \"7-4=\", \"51-13=\", \"50+10=\", \"42-39=\", \"3+76=\", \"50-18=\", \"39+14=\", \"31+63=\", \"27+17=\", \"39-31=\", \"65-1=\", \"59-3=\", \"95-42=\", \"42+5=\", \"33+34=\", \"81-49=\", \"54-46=\", \"70-62=\", \"49-47=\", \"82-5=\", \"34+41=\", \"49-42=\", \"11+56=\", \"28+31=\", \"60-56=\", \"76-18=\", \"87-65=\", \"8+17=\", \"30+15=\", \"96-46=\", \"97-57=\", \"83-75=\", \"87-23=\", \"22+2=\", \"62-10=\", \"78+1=\", \"76-45=\", \"72-56=\", \"46+33=\", \"19+24=\", \"35-3=\", \"81-49=\", \"90-70=\", \"35+59=\", \"40-9=\", \"59-48=\", \"48+34=\", \"95-69=\", \"66+11=\", \"34+41=\", \"99-76=\", \"59+5=\", \"33+26=\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = Math.min(newTexts.length, paragraphs.items.length);\nfor (let i = 0; i < count; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date line (first paragraph, not part of the table).\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = '2025-03-21 Friday'\n\n# Update every cell of the single 20x5 table, in row-major order,\n# matching the order cells appear in the document/diff.\n$t = $d.Tables.Item(1)\n$values = @(\n    @('78-78=', '45-1=', '4+38=', '3+83=', '83-42='),\n    @('72-36=', '54+30=', '37-13=', '60-29=', '46-13='),\n    @('82-1=', '91-65=', '57-22=', '75-7=', '69+10='),\n    @('42+0=', '16-12=', '39+2=', '78-49=', '77-52='),\n    @('72+10=', '47-3=', '9+6=', '31+22=', '63+6='),\n    @('60+16=', '77-25=', '78-47=', '30+28=', '24+61='),\n    @('32-20=', '28-19=', '89-64=', '57-34=', '14+72='),\n    @('23+53=', '17+16=', '50+45=', '30+46=', '17+56='),\n    @('42+15=', '13+32=', '74-37=', '20+57=', '55+39='),\n    @('45+17=', '55+19=', '7-4=', '51-13=', '50+10='),\n    @('42-39=', '3+76=', '50-18=', '39+14=', '31+63='),\n    @('27+17=', '39-31=', '65-1=', '59-3=', '95-42='),\n    @('42+5=', '33+34=', '81-49=', '54-46=', '70-62='),\n    @('49-47=', '82-5=', '34+41=', '49-42=', '11+56='),\n    @('28+31=', '60-56=', '76-18=', '87-65=', '8+17='),\n    @('30+15=', '96-46=', '97-57=', '83-75=', '87-23='),\n    @('22+2=', '62-10=', '78+1=', '76-45=', '72-56='),\n    @('46+33=', '19+24=', '35-3=', '81-49=', '90-70='),\n    @('35+59=', '40-9=', '59-48=', '48+34=', '95-69='),\n    @('66+11=', '34+41=', '99-76=', '59+5=', '33+26=')\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$r - 1][$c - 1]\n    }\n}\n"}
